$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 20) below the last existing row (row 19),
# mirroring its formatting (A column keeps the date-style format).
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.622852459381209
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 1.946625946175717
